$wb = $excel.ActiveWorkbook
$wsConst = $wb.Worksheets.Item("constants")
$wsTime = $wb.Worksheets.Item("time_variants")

# --- Insert two new rows after current row 36 (before current row 37) ---
# Copy formatting from row 36 (which already has the style used by the new rows)
$wsConst.Rows.Item(36).Copy() | Out-Null
$wsConst.Rows.Item(37).Insert() | Out-Null
$wsConst.Rows.Item(36).Copy() | Out-Null
$wsConst.Rows.Item(37).Insert() | Out-Null
$excel.CutCopyMode = $false

# New row 37: epi_prop_smearpos
$wsConst.Range("A37").Value = "epi_prop_smearpos"
$wsConst.Range("B37").Value = 0.4
$wsConst.Range("C37").Value = $null
$wsConst.Range("D37").Value = $null
$wsConst.Range("E37").Value = $null

# New row 38: epi_prop_smearneg
$wsConst.Range("A38").Value = "epi_prop_smearneg"
$wsConst.Range("B38").Value = 0.4
$wsConst.Range("C38").Value = $null
$wsConst.Range("D38").Value = $null
$wsConst.Range("E38").Value = $null

# --- Update view state ---
# constants sheet becomes the active/selected tab, with a specific scroll
# position and selection.
$wsConst.Activate()
$wsConst.Application.ActiveWindow.ScrollRow = 15
$wsConst.Range("B39").Select() | Out-Null

# time_variants should no longer be the tab that is marked selected
$wsTime.Select() | Out-Null
$wsConst.Activate()
